$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Range("E2").Value = 3
$ws.Range("G2").Value = 19.44654466666667
$ws.Range("H2").Value = 58.339634
$ws.Range("I2").Value = 0.7934109702307454
$ws.Range("J2").Value = 0.7934109702307454
$ws.Range("K2").Value = 3
$ws.Range("M2").Value = 0.4806873333333333
$ws.Range("N2").Value = 1.442062
$ws.Range("O2").Value = 0.06311654432781515
$ws.Range("P2").Value = 0.06311654432781515
$ws.Range("Q2").Value = 9.347707698367556
$ws.Range("R2").Value = 84.129369285308
$ws.Range("S2").Value = 0.05007735867274366
$ws.Range("T2").Value = 0.05007735867274366

# Row 3
$ws.Range("E3").Value = 3
$ws.Range("G3").Value = 19.44654466666667
$ws.Range("H3").Value = 58.339634
$ws.Range("I3").Value = 0.7934109702307454
$ws.Range("J3").Value = 0.7934109702307454
$ws.Range("K3").Value = 3
$ws.Range("M3").Value = 3.424957333333333
$ws.Range("N3").Value = 10.274872
$ws.Range("O3").Value = 0.44971326756452
$ws.Range("P3").Value = 0.44971326756452
$ws.Range("Q3").Value = 66.60358576409423
$ws.Range("R3").Value = 599.4322718768481
$ws.Range("S3").Value = 0.3568074399440045
$ws.Range("T3").Value = 0.3568074399440045

# Row 4
$ws.Range("E4").Value = 3
$ws.Range("G4").Value = 19.44654466666667
$ws.Range("H4").Value = 58.339634
$ws.Range("I4").Value = 0.7934109702307454
$ws.Range("J4").Value = 0.7934109702307454
$ws.Range("K4").Value = 3
$ws.Range("M4").Value = 3.710224333333334
$ws.Range("N4").Value = 11.130673
$ws.Range("O4").Value = 0.4871701881076649
$ws.Range("P4").Value = 0.4871701881076648
$ws.Range("Q4").Value = 72.15104322152024
$ws.Range("R4").Value = 649.3593889936822
$ws.Range("S4").Value = 0.3865261716139971
$ws.Range("T4").Value = 0.386526171613997

# Row 5
$ws.Range("E5").Value = 3
$ws.Range("G5").Value = 2.524415666666667
$ws.Range("H5").Value = 7.573247
$ws.Range("I5").Value = 0.1029951139231878
$ws.Range("J5").Value = 0.1029951139231878
$ws.Range("K5").Value = 3
$ws.Range("M5").Value = 0.4806873333333333
$ws.Range("N5").Value = 1.442062
$ws.Range("O5").Value = 0.06311654432781515
$ws.Range("P5").Value = 0.06311654432781515
$ws.Range("Q5").Value = 1.213454635034889
$ws.Range("R5").Value = 10.921091715314
$ws.Range("S5").Value = 0.006500695673481255
$ws.Range("T5").Value = 0.006500695673481255

# Row 6
$ws.Range("E6").Value = 3
$ws.Range("G6").Value = 2.524415666666667
$ws.Range("H6").Value = 7.573247
$ws.Range("I6").Value = 0.1029951139231878
$ws.Range("J6").Value = 0.1029951139231878
$ws.Range("K6").Value = 3
$ws.Range("M6").Value = 3.424957333333333
$ws.Range("N6").Value = 10.274872
$ws.Range("O6").Value = 0.44971326756452
$ws.Range("P6").Value = 0.44971326756452
$ws.Range("Q6").Value = 8.646015949931556
$ws.Range("R6").Value = 77.814143549384
$ws.Range("S6").Value = 0.04631826922557677
$ws.Range("T6").Value = 0.04631826922557677

# Row 7
$ws.Range("E7").Value = 3
$ws.Range("G7").Value = 2.524415666666667
$ws.Range("H7").Value = 7.573247
$ws.Range("I7").Value = 0.1029951139231878
$ws.Range("J7").Value = 0.1029951139231878
$ws.Range("K7").Value = 3
$ws.Range("M7").Value = 3.710224333333334
$ws.Range("N7").Value = 11.130673
$ws.Range("O7").Value = 0.4871701881076649
$ws.Range("P7").Value = 0.4871701881076648
$ws.Range("Q7").Value = 9.366148433914557
$ws.Range("R7").Value = 84.29533590523101
$ws.Range("S7").Value = 0.05017614902412978
$ws.Range("T7").Value = 0.05017614902412978

# Row 8
$ws.Range("E8").Value = 3
$ws.Range("G8").Value = 2.539092333333333
$ws.Range("H8").Value = 7.617277
$ws.Range("I8").Value = 0.1035939158460669
$ws.Range("J8").Value = 0.1035939158460669
$ws.Range("K8").Value = 3
$ws.Range("M8").Value = 0.4806873333333333
$ws.Range("N8").Value = 1.442062
$ws.Range("O8").Value = 0.06311654432781515
$ws.Range("P8").Value = 0.06311654432781515
$ws.Range("Q8").Value = 1.220509522797111
$ws.Range("R8").Value = 10.984585705174
$ws.Range("S8").Value = 0.00653848998159023
$ws.Range("T8").Value = 0.00653848998159023

# Row 9
$ws.Range("E9").Value = 3
$ws.Range("G9").Value = 2.539092333333333
$ws.Range("H9").Value = 7.617277
$ws.Range("I9").Value = 0.1035939158460669
$ws.Range("J9").Value = 0.1035939158460669
$ws.Range("K9").Value = 3
$ws.Range("M9").Value = 3.424957333333333
$ws.Range("N9").Value = 10.274872
$ws.Range("O9").Value = 0.44971326756452
$ws.Range("P9").Value = 0.44971326756452
$ws.Range("Q9").Value = 8.696282907060443
$ws.Range("R9").Value = 78.26654616354399
$ws.Range("S9").Value = 0.04658755839493862
$ws.Range("T9").Value = 0.04658755839493862

# Row 10
$ws.Range("E10").Value = 3
$ws.Range("G10").Value = 2.539092333333333
$ws.Range("H10").Value = 7.617277
$ws.Range("I10").Value = 0.1035939158460669
$ws.Range("J10").Value = 0.1035939158460669
$ws.Range("K10").Value = 3
$ws.Range("M10").Value = 3.710224333333334
$ws.Range("N10").Value = 11.130673
$ws.Range("O10").Value = 0.4871701881076649
$ws.Range("P10").Value = 0.4871701881076648
$ws.Range("Q10").Value = 9.420602159713445
$ws.Range("R10").Value = 84.78541943742101
$ws.Range("S10").Value = 0.05046786746953799
$ws.Range("T10").Value = 0.05046786746953798
